$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Room" column (K) with header and values
$ws.Range("K1").Value = "Room"
$ws.Range("K2").Value = "I42"
$ws.Range("K3").Value = "I42"
$ws.Range("K4").Value = "I42"

$ws.Range("K5").Select()
